$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts old C/D -> D/E), inheriting row formatting
# (row 2 style / row height) the same way Excel does for a normal column insert.
$ws.Columns("C").Insert()

# New header + query text for the case-detail query that now lives in column C.
$ws.Range("C1").Value() = "caseDetailQuery"
$query = "MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN ['caseid'] RETURN f.file_name AS ``File Name`` ,f.file_type AS ``File Type``,head(labels(parent)) AS ``Association``, f.file_description AS ``Description``,f.file_format AS Format,((f.file_size)/1024) AS Size"
$ws.Range("C2").Value() = $query

# Match column C's width to columns A/B so it renders the same as the other
# query columns (closest width obtainable through the ColumnWidth API).
$ws.Columns("C").ColumnWidth = 75

# Reproduce the author's on-screen selection/scroll state: cell C2 selected,
# sheet scrolled so column B is the left-most visible column.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C2").Select() | Out-Null
